$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Remove the stray "_GoBack" bookmark that currently sits on the
#    "Version 0.4" heading paragraph (it will be re-created further down,
#    around the word "soms" in the newly added sentence).
# ---------------------------------------------------------------------------
$goBack = $d.Bookmarks.Item("_GoBack")
$goBack.Delete()

# ---------------------------------------------------------------------------
# 2. Add a new bullet right after "Begintijd opslaan in State.Self.Starttime"
#    (and before "Stretch goals:"), carrying over the same list style /
#    numbering as its neighbours.
# ---------------------------------------------------------------------------
$anchor = $d.Content
$anchor.Find.Execute("State.Self.Starttime") | Out-Null
$sourcePara = $anchor.Paragraphs(1)

# Creates a new paragraph right after the source one, inheriting its
# paragraph style / numbering (Lijstalinea, ilvl=1, numId=30).
$sourcePara.Range.InsertParagraphAfter()
$newPara = $sourcePara.Next()

$newPara.Range.InsertAfter("Toegangspasje bemachtigen als eerste item, maakt soms van state.self.inventory geen array.")

# Underline the word "soms" inside the freshly inserted sentence.
$somsRange = $newPara.Range.Duplicate
$somsRange.Find.Execute("soms") | Out-Null
$somsRange.Font.Underline = 1

# Re-create the "_GoBack" bookmark around that same word.
$bmRange = $newPara.Range.Duplicate
$bmRange.Find.Execute("soms") | Out-Null
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
